$wb = $excel.ActiveWorkbook

function Set-EventRow {
    param($ws, $row, $a, $b, $c, $d, $e, $f, $g, $h, $i)
    $ws.Cells.Item($row, 1).Value = $a
    # Column B ("start time") is a plain text date string in the source data,
    # not a real date value -- force text format so COM does not coerce it
    # into a date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

$ws = $wb.Worksheets.Item("展览")

Set-EventRow $ws 2 "1" "2024-06-01" "苏州·环球港漫展·六一档" "相城大道1609号 苏州环球港" "2024.06.01 10:00-06.02 16:00" "52" "39" "https://show.bilibili.com/platform/detail.html?id=85850" "//i1.hdslb.com/bfs/openplatform/202405/Vorf4QMx1715584629317.jpeg"
Set-EventRow $ws 3 "2" "2024-06-08" "【会员购严选】苏州·back to ACG端阳嘉年华动漫国潮文化节" "金山南路288号 广电国际会展中心" "2024.06.08 10:00-06.09 17:00" "11511" "60" "https://show.bilibili.com/platform/detail.html?id=82233" "//i0.hdslb.com/bfs/openplatform/202405/vPI9YxcW1715674161718.jpeg"
Set-EventRow $ws 4 "3" "2024-06-15" "苏州·蔚蓝档案ONLY#2024~Game Builders Go!!!!" "城际路21号 苏州汇融广场假日酒店" "2024.06.15 10:00-06.15 17:00" "206" "75" "https://show.bilibili.com/platform/detail.html?id=84130" "//i0.hdslb.com/bfs/openplatform/202404/bpTzFcDq1713253785881.jpeg"
Set-EventRow $ws 5 "4" "2024-06-29" "苏州·归离之缘原神only展" "清禾路888号2号楼3楼 格莱美婚礼宴会中心" "2024.06.29 09:30-06.29 18:30" "325" "89" "https://show.bilibili.com/platform/detail.html?id=83271" "//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png"
Set-EventRow $ws 6 "5" "2024-07-06" "苏州·第一届寒假动漫展宅舞比赛-CF01" "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店" "2024.07.06 10:00-07.06 16:00" "221" "49" "https://show.bilibili.com/platform/detail.html?id=80528" "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"
Set-EventRow $ws 7 "6" "2024-07-19" "苏州·萤火国潮文化节动漫品牌博览会" "金山南路288号木渎影视城F2 苏州广电国际会展中心" "2024.07.19 10:00-07.21 17:00" "11467" "60" "https://show.bilibili.com/platform/detail.html?id=83301" "//i1.hdslb.com/bfs/openplatform/202405/Eh06dOvF1715926655440.jpeg"
Set-EventRow $ws 8 "7" "2024-07-20" "苏州·白日梦想7.20全职猎人ONLY展" "金芳路与新发路交叉口东南120米 万龙大厦" "2024.07.20 09:00-07.20 17:00" "469" "72" "https://show.bilibili.com/platform/detail.html?id=83508" "//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg"
Set-EventRow $ws 9 "8" "2024-07-27" "苏州·第一届动漫游戏展" "清禾路886号 尹山湖大剧院" "2024.07.27 10:30-07.27 17:00" "1161" "60" "https://show.bilibili.com/platform/detail.html?id=84899" "//i2.hdslb.com/bfs/openplatform/202404/ARz0BVLv1712661597595.jpeg"
Set-EventRow $ws 10 "9" "2024-08-03" "常熟·ACG动漫游戏嘉年华" "冬青路88号 江南·美好汇生活广场" "2024.08.03 09:00-08.04 17:00" "80" "60" "https://show.bilibili.com/platform/detail.html?id=85851" "//i2.hdslb.com/bfs/openplatform/202405/LgJRjcDn1715933608635.jpeg"
Set-EventRow $ws 11 "10" "2024-08-03" "苏州·星部落动漫嘉年华" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.08.03 09:00-08.04 16:00" "1752" "49" "https://show.bilibili.com/platform/detail.html?id=84858" "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"
Set-EventRow $ws 12 "11" "2024-08-17" "苏州·ICAN summer World动漫品牌夏游节" "金山南路288号木渎影视城F2 苏州广电国际会展中心" "2024.08.17 10:00-08.18 17:00" "5691" "60" "https://show.bilibili.com/platform/detail.html?id=85289" "//i0.hdslb.com/bfs/openplatform/202404/JavlW9fj1714459472747.jpeg"
Set-EventRow $ws 13 "12" "2024-08-17" "苏州·第二届Redamancy动漫游戏嘉年华" "清禾路886号 尹山湖大剧院" "2024.08.17 10:00-08.18 17:00" "108" "60" "https://show.bilibili.com/platform/detail.html?id=83576" "//i1.hdslb.com/bfs/openplatform/202405/UbwYg1jn1716516632746.jpeg"
Set-EventRow $ws 14 "13" "2024-10-01" "苏州·第十三届理想乡动漫展-同人创作者大会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.10.01 10:00-10.03 17:00" "3492" "39" "https://show.bilibili.com/platform/detail.html?id=83821" "//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg"
Set-EventRow $ws 15 "14" "2024-10-02" "苏州·明日方舟ONLY#2024~佑桑柔" "城际路21号 苏州汇融广场假日酒店" "2024.10.02 10:00-10.02 17:00" "178" "75" "https://show.bilibili.com/platform/detail.html?id=84046" "//i2.hdslb.com/bfs/openplatform/202405/0VhIRprD1716344515303.jpeg"
Set-EventRow $ws 16 "15" "2024-10-26" "苏州·第三届华盟国漫次元嘉年华" "清禾路886号 苏州聚橙尹山湖大剧院" "2024.10.26 10:00-10.27 17:00" "13" "58" "https://show.bilibili.com/platform/detail.html?id=85767" "//i1.hdslb.com/bfs/openplatform/202405/CqSYBZhQ1715846719965.jpeg"

$ws.Rows.Item(18).Delete()
$ws.Rows.Item(17).Delete()

$ws = $wb.Worksheets.Item("全部类型")

Set-EventRow $ws 2 "1" "2024-06-01" "苏州·春日计划2024——特别二次元不插电音乐会" "星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院" "2024.06.01 19:30-06.01 21:00" "569" "88" "https://show.bilibili.com/platform/detail.html?id=84720" "//i1.hdslb.com/bfs/openplatform/202405/MTs1Gl1Z1715588874037.jpeg"
Set-EventRow $ws 3 "2" "2024-06-01" "苏州·环球港漫展·六一档" "相城大道1609号 苏州环球港" "2024.06.01 10:00-06.02 16:00" "52" "39" "https://show.bilibili.com/platform/detail.html?id=85850" "//i1.hdslb.com/bfs/openplatform/202405/Vorf4QMx1715584629317.jpeg"
Set-EventRow $ws 4 "3" "2024-06-02" "苏州·英雄时代2024哈瓦西钢琴演奏会" "东太湖大道12000号 苏州湾大剧院" "2024.06.02 19:30-06.02 21:00" "2" "499" "https://show.bilibili.com/platform/detail.html?id=83901" "//i0.hdslb.com/bfs/openplatform/202404/LbCirky11712569675168.png"
Set-EventRow $ws 5 "4" "2024-06-08" "【会员购严选】苏州·back to ACG端阳嘉年华动漫国潮文化节" "金山南路288号 广电国际会展中心" "2024.06.08 10:00-06.09 17:00" "11511" "60" "https://show.bilibili.com/platform/detail.html?id=82233" "//i0.hdslb.com/bfs/openplatform/202405/vPI9YxcW1715674161718.jpeg"
Set-EventRow $ws 6 "5" "2024-06-15" "苏州·蔚蓝档案ONLY#2024~Game Builders Go!!!!" "城际路21号 苏州汇融广场假日酒店" "2024.06.15 10:00-06.15 17:00" "206" "75" "https://show.bilibili.com/platform/detail.html?id=84130" "//i0.hdslb.com/bfs/openplatform/202404/bpTzFcDq1713253785881.jpeg"
Set-EventRow $ws 7 "6" "2024-06-29" "苏州·归离之缘原神only展" "清禾路888号2号楼3楼 格莱美婚礼宴会中心" "2024.06.29 09:30-06.29 18:30" "325" "89" "https://show.bilibili.com/platform/detail.html?id=83271" "//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png"
Set-EventRow $ws 8 "7" "2024-07-06" "苏州·第一届寒假动漫展宅舞比赛-CF01" "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店" "2024.07.06 10:00-07.06 16:00" "221" "49" "https://show.bilibili.com/platform/detail.html?id=80528" "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"
Set-EventRow $ws 9 "8" "2024-07-19" "苏州·萤火国潮文化节动漫品牌博览会" "金山南路288号木渎影视城F2 苏州广电国际会展中心" "2024.07.19 10:00-07.21 17:00" "11467" "60" "https://show.bilibili.com/platform/detail.html?id=83301" "//i1.hdslb.com/bfs/openplatform/202405/Eh06dOvF1715926655440.jpeg"
Set-EventRow $ws 10 "9" "2024-07-20" "苏州·白日梦想7.20全职猎人ONLY展" "金芳路与新发路交叉口东南120米 万龙大厦" "2024.07.20 09:00-07.20 17:00" "469" "72" "https://show.bilibili.com/platform/detail.html?id=83508" "//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg"
Set-EventRow $ws 11 "10" "2024-07-27" "苏州·第一届动漫游戏展" "清禾路886号 尹山湖大剧院" "2024.07.27 10:30-07.27 17:00" "1161" "60" "https://show.bilibili.com/platform/detail.html?id=84899" "//i2.hdslb.com/bfs/openplatform/202404/ARz0BVLv1712661597595.jpeg"
Set-EventRow $ws 12 "11" "2024-08-03" "常熟·ACG动漫游戏嘉年华" "冬青路88号 江南·美好汇生活广场" "2024.08.03 09:00-08.04 17:00" "80" "60" "https://show.bilibili.com/platform/detail.html?id=85851" "//i2.hdslb.com/bfs/openplatform/202405/LgJRjcDn1715933608635.jpeg"
Set-EventRow $ws 13 "12" "2024-08-03" "苏州·星部落动漫嘉年华" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.08.03 09:00-08.04 16:00" "1752" "49" "https://show.bilibili.com/platform/detail.html?id=84858" "//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg"
Set-EventRow $ws 14 "13" "2024-08-10" "苏州·爱乐之城·经典电影作品音乐会" "念珠街121号道前街与吉庆街路口距养育巷地铁站 苏州市会议中心" "2024.08.10 19:30-08.10 21:00" "2" "50" "https://show.bilibili.com/platform/detail.html?id=86194" "//i2.hdslb.com/bfs/openplatform/202405/vagzbfox1716438290025.jpeg"
Set-EventRow $ws 15 "14" "2024-08-17" "苏州·ICAN summer World动漫品牌夏游节" "金山南路288号木渎影视城F2 苏州广电国际会展中心" "2024.08.17 10:00-08.18 17:00" "5691" "60" "https://show.bilibili.com/platform/detail.html?id=85289" "//i0.hdslb.com/bfs/openplatform/202404/JavlW9fj1714459472747.jpeg"
Set-EventRow $ws 16 "15" "2024-08-17" "苏州·第二届Redamancy动漫游戏嘉年华" "清禾路886号 尹山湖大剧院" "2024.08.17 10:00-08.18 17:00" "108" "60" "https://show.bilibili.com/platform/detail.html?id=83576" "//i1.hdslb.com/bfs/openplatform/202405/UbwYg1jn1716516632746.jpeg"
Set-EventRow $ws 17 "16" "2024-10-01" "苏州·第十三届理想乡动漫展-同人创作者大会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.10.01 10:00-10.03 17:00" "3492" "39" "https://show.bilibili.com/platform/detail.html?id=83821" "//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg"
Set-EventRow $ws 18 "17" "2024-10-02" "苏州·明日方舟ONLY#2024~佑桑柔" "城际路21号 苏州汇融广场假日酒店" "2024.10.02 10:00-10.02 17:00" "178" "75" "https://show.bilibili.com/platform/detail.html?id=84046" "//i2.hdslb.com/bfs/openplatform/202405/0VhIRprD1716344515303.jpeg"
Set-EventRow $ws 19 "18" "2024-10-26" "苏州·第三届华盟国漫次元嘉年华" "清禾路886号 苏州聚橙尹山湖大剧院" "2024.10.26 10:00-10.27 17:00" "13" "58" "https://show.bilibili.com/platform/detail.html?id=85767" "//i1.hdslb.com/bfs/openplatform/202405/CqSYBZhQ1715846719965.jpeg"

$ws.Rows.Item(21).Delete()
$ws.Rows.Item(20).Delete()

